$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the authoritative diff.
$updates = [ordered]@{
    'D2' = '29.624.43'
    'E2' = '  +2.54%  '
    'D3' = '1.861.01'
    'D4' = '0.9992'
    'E4' = '  -0.01%  '
    'D5' = '245.28'
    'E5' = '  +1.68%  '
    'D6' = '0.6985'
    'E6' = '  +1.19%  '
    'D7' = '0.9998'
    'E7' = '  +0.03%  '
    'E8' = '  +1.19%  '
    'D9' = '0.3064'
    'E9' = '  +0.53%  '
    'D10' = '23.69'
    'E10' = '  +0.91%  '
    'D11' = '0.07752'
    'E11' = '  -0.14%  '
    'D12' = '5.163'
    'E12' = '  +2.27%  '
    'D13' = '1.859.45'
    'E13' = '  +1.77%  '
    'D14' = '92.38'
    'E14' = '  +2.30%  '
    'D15' = '0.6922'
    'E15' = '  +2.85%  '
    'D16' = '6.571'
    'E16' = '  +3.04%  '
    'D17' = '29.599.65'
    'E17' = '  +2.52%  '
    'D18' = '0.000008334'
    'E18' = '  +0.79%  '
    'D19' = '2.105.29'
    'E19' = '  +1.60%  '
    'D20' = '242.06'
    'E20' = '  +0.03%  '
    'D21' = '12.77'
    'D22' = '0.9998'
    'E22' = '  +0.03%  '
    'D23' = '7.620'
    'E23' = '  +2.81%  '
    'E24' = '  +0.11%  '
    'D25' = '0.1503'
    'E25' = '  +2.40%  '
    'D26' = '8.923'
    'E26' = '  +2.05%  '
    'D27' = '159.39'
    'E27' = '  -0.90%  '
    'E28' = '  +0.70%  '
    'D29' = '1.534'
    'E29' = '  -0.19%  '
    'D30' = '4.255'
    'D31' = '4.187'
    'E31' = '  +1.31%  '
    'D32' = '1.195'
    'E32' = '  -0.01%  '
    'D33' = '0.05091'
    'E33' = '  -0.13%  '
    'D34' = '0.7769'
    'E34' = '  +4.07%  '
    'D35' = '1.897'
    'E35' = '  +4.75%  '
    'D36' = '1.155'
    'D37' = '2.684'
    'E37' = '  +0.30%  '
    'D38' = '1.326.66'
    'E38' = '  +10.27%  '
    'D39' = '0.01874'
    'E39' = '  +2.02%  '
    'D40' = '2.731'
    'E40' = '  +2.14%  '
    'D41' = '0.9624'
    'E41' = '  +3.86%  '
    'D42' = '106.50'
    'E42' = '  -1.64%  '
    'D43' = '5.818'
    'E43' = '  +11.84%  '
    'D44' = '0.9996'
    'E44' = '  +0.03%  '
    'D45' = '9.774'
    'E45' = '  +3.30%  '
    'B46' = 'RocketPoolETH'
    'C46' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D46' = '2.005.12'
    'E46' = '  +1.58%  '
    'B47' = 'BabyDogeCoin'
    'C47' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D47' = '0.00000000124'
    'E47' = '  +2.95%  '
    'D48' = '0.5216'
    'E48' = '  +1.04%  '
    'B49' = 'RenderToken'
    'C49' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D49' = '1.781'
    'E49' = '  +3.29%  '
    'B50' = 'Aave'
    'C50' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D50' = '64.38'
    'E50' = '  +4.07%  '
    'D51' = '6.977'
    'E51' = '  +1.49%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $newVal = $updates[$addr]
    if ($addr[0] -eq "D") {
        # Column D ("Price") mixes genuine numeric-looking text (e.g. "0.9992")
        # with thousand-dot-separated text that is not a valid number
        # (e.g. "29.624.43"). The source file stores ALL of these as plain
        # text, so force text formatting before writing, then drop the
        # explicit NumberFormat again so the cell keeps using the default
        # (unstyled) cell style - only the stored value/type should change.
        $cell.NumberFormat = "@"
        $cell.Value = $newVal
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newVal
    }
}
